# Update simulated game probability matrix values on Sheet1
# (added more games, sped up simulate game logic, drafted optimization logic)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1962025316455696
$ws.Range("C2").Value = 0.560126582278481
$ws.Range("J2").Value = 0.0189873417721519
$ws.Range("P2").Value = 0.1424050632911392
$ws.Range("S2").Value = 0.08227848101265822
$ws.Range("B3").Value = 0.01630434782608696
$ws.Range("C3").Value = 0.02717391304347826
$ws.Range("J3").Value = 0.03260869565217391
$ws.Range("P3").Value = 0.7608695652173914
$ws.Range("S3").Value = 0.1630434782608696
$ws.Range("J4").Value = 0.02173913043478261
$ws.Range("P4").Value = 0.7608695652173914
$ws.Range("S4").Value = 0.2173913043478261
$ws.Range("B6").Value = 0.07727272727272727
$ws.Range("D6").Value = 0.01818181818181818
$ws.Range("F6").Value = 0.07727272727272727
$ws.Range("J6").Value = 0.2863636363636364
$ws.Range("O6").Value = 0.01363636363636364
$ws.Range("Q6").Value = 0.1954545454545455
$ws.Range("R6").Value = 0.05909090909090909
$ws.Range("S6").Value = 0.2727272727272727
$ws.Range("B7").Value = 0.1355311355311355
$ws.Range("D7").Value = 0.01098901098901099
$ws.Range("F7").Value = 0.0293040293040293
$ws.Range("J7").Value = 0.1428571428571428
$ws.Range("O7").Value = 0.01831501831501832
$ws.Range("Q7").Value = 0.2051282051282051
$ws.Range("R7").Value = 0.09523809523809523
$ws.Range("S7").Value = 0.3626373626373626
$ws.Range("B8").Value = 0.07312252964426877
$ws.Range("D8").Value = 0.02569169960474308
$ws.Range("F8").Value = 0.04347826086956522
$ws.Range("J8").Value = 0.1146245059288538
$ws.Range("O8").Value = 0.01185770750988142
$ws.Range("Q8").Value = 0.2391304347826087
$ws.Range("R8").Value = 0.09486166007905138
$ws.Range("S8").Value = 0.3972332015810277
$ws.Range("B9").Value = 0.09663865546218488
$ws.Range("D9").Value = 0.02100840336134454
$ws.Range("F9").Value = 0.06722689075630252
$ws.Range("J9").Value = 0.07563025210084033
$ws.Range("O9").Value = 0.01260504201680672
$ws.Range("Q9").Value = 0.180672268907563
$ws.Range("R9").Value = 0.07563025210084033
$ws.Range("S9").Value = 0.4705882352941176
$ws.Range("B10").Value = 0.1032403918613414
$ws.Range("D10").Value = 0.01582516955538809
$ws.Range("F10").Value = 0.06631499623210249
$ws.Range("J10").Value = 0.132629992464205
$ws.Range("O10").Value = 0.02260738507912585
$ws.Range("Q10").Value = 0.2426525998492841
$ws.Range("R10").Value = 0.08741522230595328
$ws.Range("S10").Value = 0.3293142426525998
$ws.Range("G11").Value = 0.1620253164556962
$ws.Range("J11").Value = 0.06835443037974684
$ws.Range("K11").Value = 0.1873417721518987
$ws.Range("L11").Value = 0.569620253164557
$ws.Range("S11").Value = 0.01265822784810127
$ws.Range("F12").Value = 0.00423728813559322
$ws.Range("G12").Value = 0.7372881355932204
$ws.Range("J12").Value = 0.173728813559322
$ws.Range("K12").Value = 0.0211864406779661
$ws.Range("L12").Value = 0.03813559322033899
$ws.Range("S12").Value = 0.02542372881355932
$ws.Range("G13").Value = 0.6825396825396826
$ws.Range("J13").Value = 0.2380952380952381
$ws.Range("S13").Value = 0.07936507936507936
$ws.Range("F15").Value = 0.02714932126696833
$ws.Range("H15").Value = 0.1493212669683258
$ws.Range("I15").Value = 0.04072398190045249
$ws.Range("J15").Value = 0.3393665158371041
$ws.Range("K15").Value = 0.08597285067873303
$ws.Range("M15").Value = 0.01809954751131222
$ws.Range("O15").Value = 0.06787330316742081
$ws.Range("S15").Value = 0.2714932126696832
$ws.Range("F16").Value = 0.01408450704225352
$ws.Range("H16").Value = 0.2065727699530517
$ws.Range("I16").Value = 0.1032863849765258
$ws.Range("J16").Value = 0.3661971830985916
$ws.Range("K16").Value = 0.1267605633802817
$ws.Range("M16").Value = 0.02816901408450704
$ws.Range("N16").Value = 0.009389671361502348
$ws.Range("O16").Value = 0.03286384976525822
$ws.Range("S16").Value = 0.1126760563380282
$ws.Range("F17").Value = 0.01896551724137931
$ws.Range("H17").Value = 0.1931034482758621
$ws.Range("I17").Value = 0.09827586206896552
$ws.Range("J17").Value = 0.4
$ws.Range("K17").Value = 0.1120689655172414
$ws.Range("M17").Value = 0.02068965517241379
$ws.Range("O17").Value = 0.05689655172413793
$ws.Range("S17").Value = 0.1
$ws.Range("F18").Value = 0.004524886877828055
$ws.Range("H18").Value = 0.1945701357466063
$ws.Range("I18").Value = 0.09502262443438914
$ws.Range("J18").Value = 0.4208144796380091
$ws.Range("K18").Value = 0.06334841628959276
$ws.Range("M18").Value = 0.02262443438914027
$ws.Range("N18").Value = 0.004524886877828055
$ws.Range("O18").Value = 0.05429864253393665
$ws.Range("S18").Value = 0.1402714932126697
$ws.Range("F19").Value = 0.01379310344827586
$ws.Range("H19").Value = 0.2068965517241379
$ws.Range("I19").Value = 0.09885057471264368
$ws.Range("J19").Value = 0.3371647509578544
$ws.Range("K19").Value = 0.1363984674329502
$ws.Range("M19").Value = 0.02605363984674329
$ws.Range("N19").Value = 0.0007662835249042146
$ws.Range("O19").Value = 0.06130268199233716
$ws.Range("S19").Value = 0.1187739463601533
